$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 613.3333
$ws.Range("I41").Value = 507.07693
$ws.Range("J41").Value = 786
$ws.Range("K41").Value = 507.07693
$ws.Range("L41").Value = 786
$ws.Range("M41").Value = -67.07693
$ws.Range("N41").Value = -1666
$ws.Range("H80").Value = 8255.111000000001
$ws.Range("I80").Value = 1758.9
$ws.Range("J80").Value = 16375.375
$ws.Range("K80").Value = 5276.700000000001
$ws.Range("L80").Value = 49126.125
$ws.Range("M80").Value = -4278.700000000001
$ws.Range("N80").Value = -51122.125
$ws.Range("H83").Value = 8255.111000000001
$ws.Range("I83").Value = 1758.9
$ws.Range("J83").Value = 16375.375
$ws.Range("K83").Value = 15830.1
$ws.Range("L83").Value = 147378.375
$ws.Range("M83").Value = -10838.1
$ws.Range("N83").Value = -157362.375
$ws.Range("H98").Value = 1670.64
$ws.Range("I98").Value = 842.25
$ws.Range("K98").Value = 842.25
$ws.Range("M98").Value = 655.75
$ws.Range("H110").Value = 67899.836
$ws.Range("J110").Value = 67899.836
$ws.Range("L110").Value = 67899.836
$ws.Range("N110").Value = -76079.836
$ws.Range("H122").Value = 1670.64
$ws.Range("I122").Value = 842.25
$ws.Range("K122").Value = 2526.75
$ws.Range("M122").Value = -76.75
$ws.Range("H123").Value = 68078.664
$ws.Range("J123").Value = 68078.664
$ws.Range("L123").Value = 68078.664
$ws.Range("N123").Value = -77878.664
$ws.Range("H134").Value = 67146.625
$ws.Range("J134").Value = 67146.625
$ws.Range("L134").Value = 67146.625
$ws.Range("N134").Value = -77286.625
$ws.Range("H137").Value = 209206.58
$ws.Range("I137").Value = 1880.9783
$ws.Range("K137").Value = 5642.9349
$ws.Range("M137").Value = -3092.9349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 52730.668
$ws.Range("J52").Value = 52730.668
$ws.Range("L52").Value = 52730.668
$ws.Range("N52").Value = -53366.668
$ws.Range("H61").Value = 1840.6923
$ws.Range("I61").Value = 1618.4
$ws.Range("K61").Value = 1618.4
$ws.Range("M61").Value = -1406.4
$ws.Range("H74").Value = 2656.2415
$ws.Range("I74").Value = 1707.2778
$ws.Range("K74").Value = 1707.2778
$ws.Range("M74").Value = -833.2778000000001
$ws.Range("H77").Value = 2656.2415
$ws.Range("I77").Value = 1707.2778
$ws.Range("K77").Value = 8536.389000000001
$ws.Range("M77").Value = -4168.389000000001
$ws.Range("H104").Value = 39797.6
$ws.Range("J104").Value = 39797.6
$ws.Range("L104").Value = 39797.6
$ws.Range("N104").Value = -46785.6
$ws.Range("H110").Value = 930.6
$ws.Range("I110").Value = 811.7778
$ws.Range("K110").Value = 811.7778
$ws.Range("M110").Value = 1233.2222
$ws.Range("H122").Value = 2757.8333
$ws.Range("I122").Value = 2120.2
$ws.Range("K122").Value = 6360.599999999999
$ws.Range("M122").Value = -3910.599999999999
$ws.Range("H132").Value = 2252.1853
$ws.Range("I132").Value = 1531.2106
$ws.Range("J132").Value = 3964.5
$ws.Range("K132").Value = 4593.6318
$ws.Range("L132").Value = 11893.5
$ws.Range("M132").Value = -2063.6318
$ws.Range("N132").Value = -16953.5
$ws.Range("H136").Value = 1840.6923
$ws.Range("I136").Value = 1618.4
$ws.Range("K136").Value = 4855.200000000001
$ws.Range("M136").Value = -2305.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 15584709
$ws.Range("I22").Value = 19480836
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 19480836
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -19480663
$ws.Range("N22").Value = -546
$ws.Range("H52").Value = 99988
$ws.Range("J52").Value = 99988
$ws.Range("L52").Value = 99988
$ws.Range("N52").Value = -100514
$ws.Range("H86").Value = 5055.56
$ws.Range("I86").Value = 3906.125
$ws.Range("J86").Value = 5596.4707
$ws.Range("K86").Value = 3906.125
$ws.Range("L86").Value = 5596.4707
$ws.Range("M86").Value = -2783.125
$ws.Range("N86").Value = -7842.4707
$ws.Range("H89").Value = 5055.56
$ws.Range("I89").Value = 3906.125
$ws.Range("J89").Value = 5596.4707
$ws.Range("K89").Value = 19530.625
$ws.Range("L89").Value = 27982.3535
$ws.Range("M89").Value = -13914.625
$ws.Range("N89").Value = -39214.3535
$ws.Range("H105").Value = 26164.244
$ws.Range("I105").Value = 39691.848
$ws.Range("J105").Value = 2716.4
$ws.Range("K105").Value = 39691.848
$ws.Range("L105").Value = 2716.4
$ws.Range("M105").Value = -37944.848
$ws.Range("N105").Value = -6210.4
$ws.Range("H109").Value = 95991
$ws.Range("J109").Value = 95991
$ws.Range("L109").Value = 95991
$ws.Range("N109").Value = -98765
$ws.Range("H119").Value = 49997
$ws.Range("J119").Value = 49997
$ws.Range("L119").Value = 49997
$ws.Range("N119").Value = -59673
$ws.Range("H121").Value = 99988
$ws.Range("J121").Value = 99988
$ws.Range("L121").Value = 99988
$ws.Range("N121").Value = -103482
$ws.Range("H132").Value = 92354.28999999999
$ws.Range("J132").Value = 92354.28999999999
$ws.Range("L132").Value = 92354.28999999999
$ws.Range("N132").Value = -102474.29

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5872.25
$ws.Range("I86").Value = 5079.4
$ws.Range("K86").Value = 5079.4
$ws.Range("M86").Value = -3956.4
$ws.Range("H89").Value = 5872.25
$ws.Range("I89").Value = 5079.4
$ws.Range("K89").Value = 25397
$ws.Range("M89").Value = -19781
$ws.Range("H94").Value = 884.5
$ws.Range("I94").Value = 776
$ws.Range("J94").Value = 993
$ws.Range("K94").Value = 776
$ws.Range("L94").Value = 993
$ws.Range("M94").Value = -325
$ws.Range("N94").Value = -1895
$ws.Range("H108").Value = 57681.453
$ws.Range("J108").Value = 57681.453
$ws.Range("L108").Value = 57681.453
$ws.Range("N108").Value = -65361.453
$ws.Range("H119").Value = 63404.285
$ws.Range("J119").Value = 63404.285
$ws.Range("L119").Value = 63404.285
$ws.Range("N119").Value = -73080.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 24999.5
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 24999.5
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 69846.55499999999
$ws.Range("J109").Value = 69846.55499999999
$ws.Range("L109").Value = 69846.55499999999
$ws.Range("N109").Value = -71926.55499999999
$ws.Range("H110").Value = 74389.5
$ws.Range("J110").Value = 74389.5
$ws.Range("L110").Value = 74389.5
$ws.Range("N110").Value = -82569.5
$ws.Range("H114").Value = 91940.28999999999
$ws.Range("J114").Value = 91940.28999999999
$ws.Range("L114").Value = 91940.28999999999
$ws.Range("N114").Value = -100618.29
$ws.Range("H123").Value = 24432.75
$ws.Range("J123").Value = 24432.75
$ws.Range("L123").Value = 24432.75
$ws.Range("N123").Value = -29332.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7941936.5
$ws.Range("I40").Value = 6000.2
$ws.Range("K40").Value = 6000.2
$ws.Range("M40").Value = -5864.2
$ws.Range("H118").Value = 53700.363
$ws.Range("J118").Value = 53700.363
$ws.Range("L118").Value = 53700.363
$ws.Range("N118").Value = -57014.363
$ws.Range("H122").Value = 16671807
$ws.Range("I122").Value = 5618.1
$ws.Range("K122").Value = 16854.3
$ws.Range("M122").Value = -14404.3
$ws.Range("H123").Value = 79097.60000000001
$ws.Range("J123").Value = 79097.60000000001
$ws.Range("L123").Value = 79097.60000000001
$ws.Range("N123").Value = -88897.60000000001
$ws.Range("H136").Value = 2428.625
$ws.Range("I136").Value = 2136.5
$ws.Range("K136").Value = 6409.5
$ws.Range("M136").Value = -3859.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3777
$ws.Range("I122").Value = 2831.3333
$ws.Range("K122").Value = 8493.999899999999
$ws.Range("M122").Value = -6043.999899999999
$ws.Range("H127").Value = 93397.664
$ws.Range("J127").Value = 99999.2
$ws.Range("L127").Value = 99999.2
$ws.Range("N127").Value = -109919.2
$ws.Range("H132").Value = 1280925.9
$ws.Range("I132").Value = 1884.0625
$ws.Range("K132").Value = 5652.1875
$ws.Range("M132").Value = -3122.1875
$ws.Range("H135").Value = 61724.7
$ws.Range("J135").Value = 62280.875
$ws.Range("L135").Value = 62280.875
$ws.Range("N135").Value = -72420.875
